$wb = $excel.ActiveWorkbook

# Add the new "Booking" worksheet after the last existing sheet
# (CompoundAccessCode), so it lands at the very end of the tab strip.
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Booking"

# Header row
$ws.Range("A1").Value = "RatePlan"
$ws.Range("B1").Value = "AccessCode"
$ws.Range("C1").Value = "ArrivalD"
$ws.Range("D1").Value = "DepartureD"
$ws.Range("E1").Value = "AccessCode2"

# Data row
$ws.Range("A2").Value = "Summer Special -SS"
$ws.Range("B2").Value = "AC1"
$ws.Range("C2").Value = "'02/10/2022"
$ws.Range("D2").Value = "'02/12/2022"
$ws.Range("E2").Value = "AC1Edit"

# Column widths to match the authored layout
$ws.Columns.Item(1).ColumnWidth = 15.65
$ws.Columns.Item(2).ColumnWidth = 10.8
$ws.Columns.Item(3).ColumnWidth = 10.15
$ws.Columns.Item(4).ColumnWidth = 9.65
$ws.Columns.Item(5).ColumnWidth = 11.5

# Match the authored selection/active-cell on the new (now active) sheet
[void]$ws.Range("C4").Select()
